$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Daylight savings time update: event names no longer carry a fixed time suffix.
$ws.Range("A2").Value = "Raid the Vault 1"
$ws.Range("A3").Value = "Raid the Vault 2"
$ws.Range("A4").Value = "Shadow Assembly"
$ws.Range("A5").Value = "Battlegrounds 1"
$ws.Range("A6").Value = "Battlegrounds 2"
$ws.Range("A7").Value = "Battlegrounds 3"
$ws.Range("A8").Value = "Battlegrounds 4"
$ws.Range("A9").Value = "Shadow Lottery 1"
$ws.Range("A10").Value = "Shadow Lottery 2"
$ws.Range("A11").Value = "Shadow Lottery 3"
$ws.Range("A12").Value = "Shadow War"
$ws.Range("A13").Value = "Rite of Exile"

# Update current selection to reflect where the user left off editing.
$ws.Range("C8").Select()
